$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 21:52"

# Row 21 is now Asturias (updated figures), row 22 is now Sevilla (figures
# formerly shown on row 21) - the province list is sorted by total cases
# descending, and Asturias' updated total (1679) now outranks Sevilla's
# total (1663), so the two rows swap places.
$ws.Range("A21").Value = "Asturias"
$ws.Range("B21").Value = 1679
$ws.Range("C21").Value = 244
$ws.Range("D21").Value = 1339
$ws.Range("E21").Value = 96

$ws.Range("A22").Value = "Sevilla"
$ws.Range("B22").Value = 1663
$ws.Range("C22").Value = 95
$ws.Range("D22").Value = 1471
$ws.Range("E22").Value = 97

# Melilla (row 54) updated figures
$ws.Range("B54").Value = 92
$ws.Range("D54").Value = 79

# Ceuta (row 55) updated figures
$ws.Range("D55").Value = 68
$ws.Range("E55").Value = 4
